# CoilSetup_ProtoMPEX.xlsx -- "corrected datum column in CoilSetup_ProtoMPEX.xlsx in Example 1"
#
# Sets the "datum" column (column C) to coil index 5 on every conf_* sheet
# except conf_E, and restores each sheet's last-used selection / the
# workbook's active-tab bookmark to match the edited file.

$wb = $excel.ActiveWorkbook

# --- conf_A (sheet 1): datum column C2:C13 -> 5 -------------------------
$wsA = $wb.Worksheets.Item("conf_A")
$wsA.Range("C2:C13").Value = 5
$wsA.Range("C1:C1048576").Select()

# --- conf_B (sheet 2): datum column C2:C13 -> 5 -------------------------
$wsB = $wb.Worksheets.Item("conf_B")
$wsB.Range("C2:C13").Value = 5
$wsB.Range("E17").Select()

# --- conf_C (sheet 3): datum column C2:C14 -> 5 -------------------------
$wsC = $wb.Worksheets.Item("conf_C")
$wsC.Range("C2:C14").Value = 5
$wsC.Range("C14").Select()

# --- conf_D (sheet 4): datum column C2:C14 -> 5 (selection untouched) ---
$wsD = $wb.Worksheets.Item("conf_D")
$wsD.Range("C2:C14").Value = 5

# --- conf_E (sheet 5): left completely untouched -------------------------

# --- conf_G (sheet 7): datum column C2:C14 -> 5 --------------------------
$wsG = $wb.Worksheets.Item("conf_G")
$wsG.Range("C2:C14").Value = 5
$wsG.Range("E17").Select()

# --- conf_F (sheet 6): datum column C2:C14 -> 5 ---------------------------
# Selected/activated last so it ends up as the workbook's active tab,
# matching the saved bookViews/workbookView@activeTab bookmark.
$wsF = $wb.Worksheets.Item("conf_F")
$wsF.Range("C2:C14").Value = 5
$wsF.Activate()
$wsF.Range("F18").Select()
